# Exit info list.xlsx - edit script
# - Fill in Sheet2 ("refnum name" reference table) with the full set of rows
# - Adjust column widths on Sheet2
# - Turn off right-to-left view on Sheet2
# - Make Sheet2 the active/selected tab, with a B8 selection
# - Change Sheet1's selection to B9 (no longer the active tab)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Populate Sheet2 ---
# Row 2 - header
$ws2.Range("A2").Value = "refnum name"
$ws2.Range("B2").Value = "Device"
$ws2.Range("C2").Value = "Function "
$ws2.Range("D2").Value = "Slot"

# Row 4
$ws2.Range("A4").Value = "Drive"
$ws2.Range("B4").Value = "4 Motors"
$ws2.Range("C4").Value = "drive"
$ws2.Range("D4").Value = "PWM 1-4"

# Row 5
$ws2.Range("A5").Value = "Pilot"
$ws2.Range("B5").Value = "Joystick"
$ws2.Range("D5").Value = "USB 1"

# Row 6
$ws2.Range("A6").Value = "CoPilot"
$ws2.Range("B6").Value = "Joystick"
$ws2.Range("D6").Value = "USB 2"

# Rows 8-10
$ws2.Range("A8").Value = "First Cannon"
$ws2.Range("A9").Value = "Second Cannon"
$ws2.Range("A10").Value = "Disc Out"

# Rows 12-13
$ws2.Range("A12").Value = "CannonAI"
$ws2.Range("A13").Value = "Enc"

# Rows 15-17
$ws2.Range("A15").Value = "Flipper"
$ws2.Range("A16").Value = "Conveyer Belt"
$ws2.Range("A17").Value = "Roller"

# Rows 19-24 (row 22 entered before row 21, matching original authoring order)
$ws2.Range("A19").Value = "Cannon"
$ws2.Range("A20").Value = "InFlipper"
$ws2.Range("A22").Value = "Middle"
$ws2.Range("A21").Value = "FlipperAtPlace"
$ws2.Range("A23").Value = "Bottom"
$ws2.Range("A24").Value = "FrisbeeDirection"

# Fill in the joystick descriptions last
$ws2.Range("C5").Value = "driver's joystick"
$ws2.Range("C6").Value = "operator's joystick"

# --- Column widths on Sheet2 ---
# Target widths stored in the OOXML are 14.375 / 17.5 / 15.5 "characters".
# The engine quantizes ColumnWidth to whole pixels (MaxDigitWidth=7) when it
# serializes <col width=.../>, so these inputs land on the closest
# representable stored width (14.4286 / 17.4286 / 15.4286).
$ws2.Columns.Item(1).ColumnWidth = 13.714285714285714
$ws2.Columns.Item(2).ColumnWidth = 16.714285714285715
$ws2.Columns.Item(3).ColumnWidth = 14.714285714285714

# --- View / selection changes ---
$ws2.Activate() | Out-Null
$excel.ActiveWindow.DisplayRightToLeft = $false
$ws2.Range("B8").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("B9").Select() | Out-Null

$ws2.Activate() | Out-Null
